# Update the "last modified" date on the About sheet (C1) and
# set the Maximum Capacity Factor values to 1 for the relevant plant types
# on the MCF sheet.

$wb = $excel.ActiveWorkbook

# --- About sheet: update date value in C1 ---
$about = $wb.Worksheets.Item("About")
$about.Range("C1").Value = 45392

# --- MCF sheet: update capacity factor values ---
$mcf = $wb.Worksheets.Item("MCF")

$mcf.Range("B2").Value = 1
$mcf.Range("B3").Value = 1
$mcf.Range("B4").Value = 1
$mcf.Range("B6").Value = 1
$mcf.Range("B10").Value = 1
$mcf.Range("B11").Value = 1
$mcf.Range("B12").Value = 1
$mcf.Range("B13").Value = 1
$mcf.Range("B14").Value = 1
$mcf.Range("B16").Value = 1
$mcf.Range("B17").Value = 1
$mcf.Range("B18").Value = 1

# Select B17 and activate the MCF sheet, matching the updated cursor position
$mcf.Activate()
$mcf.Range("B17").Select()
